$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = 7.7
$ws.Range("B17").Value = 8.1
$ws.Range("B18").Value = 9.1
$ws.Range("B21").Value = 8.6
$ws.Range("B22").Value = 10.3
